# Correccion error TASK Equipo
#
# The running-total column K in the "Hoja1" tracking table accumulates
# row-by-row via K(n) = J(n) + K(n-1). Row 85 had been typed/pasted with a
# stale reference two rows back (=J85+K83) instead of the immediately
# preceding row (=J85+K84), which silently undercounted every subsequent
# cumulative total in the chain (K86:K114).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Hoja1" is already the active/selected tab
$ws.Activate()

$ws.Range("K85").Formula = "=J85+K84"

# Leave the selection where the fix was made, matching the saved view state.
$ws.Range("K84:K85").Select()
